$wb = $excel.ActiveWorkbook

# ALC!row33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 946.5333000000001
$ws.Range("I33").Value = 732.44446
$ws.Range("J33").Value = 1267.6666
$ws.Range("K33").Value = 732.44446
$ws.Range("L33").Value = 1267.6666
$ws.Range("M33").Value = -503.44446
$ws.Range("N33").Value = -1725.6666

# ALC!row125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1844.7059
$ws.Range("I125").Value = 1272.5454
$ws.Range("J125").Value = 2893.6667
$ws.Range("K125").Value = 11452.9086
$ws.Range("L125").Value = 26043.0003
$ws.Range("M125").Value = -8992.908599999999
$ws.Range("N125").Value = -30963.0003

# ALC!row129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 3847612.5
$ws.Range("J129").Value = 1613.9434
$ws.Range("L129").Value = 4841.8302
$ws.Range("N129").Value = -14841.8302

# ALC!row135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 558.73334
$ws.Range("I135").Value = 567.65515
$ws.Range("K135").Value = 5108.896350000001
$ws.Range("M135").Value = -2573.896350000001

# ALC!row138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4588.8276
$ws.Range("I138").Value = 2652.1538
$ws.Range("K138").Value = 7956.4614
$ws.Range("M138").Value = -2816.4614

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 569.3333
$ws.Range("I94").Value = 540.5
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 540.5
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = -89.5
$ws.Range("N94").Value = -1702

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2275896.8
$ws.Range("I31").Value = 2859298.8
$ws.Range("J31").Value = 7111.1113
$ws.Range("K31").Value = 2859298.8
$ws.Range("L31").Value = 7111.1113
$ws.Range("M31").Value = -2859003.8
$ws.Range("N31").Value = -7701.1113

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2275896.8
$ws.Range("I34").Value = 2859298.8
$ws.Range("J34").Value = 7111.1113
$ws.Range("K34").Value = 2859298.8
$ws.Range("L34").Value = 7111.1113
$ws.Range("M34").Value = -2859096.8
$ws.Range("N34").Value = -7515.1113

# CRP!row122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1542.1111
$ws.Range("I122").Value = 1358.35
$ws.Range("K122").Value = 4075.05
$ws.Range("M122").Value = -1625.05

# CUL!row34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 17816.666
$ws.Range("J34").Value = 17816.666
$ws.Range("L34").Value = 53449.99800000001
$ws.Range("N34").Value = -53617.99800000001

# CUL!row39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 999
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

# CUL!row47
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1721.9166
$ws.Range("I47").Value = 178.16667
$ws.Range("J47").Value = 3265.6667
$ws.Range("K47").Value = 534.50001
$ws.Range("L47").Value = 9797.000100000001
$ws.Range("M47").Value = -103.50001
$ws.Range("N47").Value = -10659.0001

# CUL!row55
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 2817.2727
$ws.Range("I55").Value = 500
$ws.Range("J55").Value = 3686.25
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 11058.75
$ws.Range("M55").Value = -1323
$ws.Range("N55").Value = -11412.75

# CUL!row82
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 313
$ws.Range("I82").Value = 313
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 939
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -533
$ws.Range("N82").ClearContents()

# CUL!row85
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H85").Value = 313
$ws.Range("I85").Value = 313
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 939
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 465
$ws.Range("N85").ClearContents()

# CUL!row87
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 8633.333000000001
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# CUL!row90
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 8633.333000000001
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# CUL!row98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 203.2
$ws.Range("I98").Value = 203.2
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 609.5999999999999
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 888.4000000000001
$ws.Range("N98").ClearContents()

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 621.069
$ws.Range("I113").Value = 499.2903
$ws.Range("J113").Value = 760.8889
$ws.Range("K113").Value = 1497.8709
$ws.Range("L113").Value = 2282.6667
$ws.Range("M113").Value = 672.1291000000001
$ws.Range("N113").Value = -6622.6667

# CUL!row131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1270.3556
$ws.Range("J131").Value = 1159.6471
$ws.Range("L131").Value = 3478.9413
$ws.Range("N131").Value = -13558.9413

# CUL!row137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2576.6
$ws.Range("J137").Value = 3665
$ws.Range("L137").Value = 10995
$ws.Range("N137").Value = -21195

# GSM!row22
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 43006.5
$ws.Range("I22").Value = 16004
$ws.Range("J22").Value = 70009
$ws.Range("K22").Value = 16004
$ws.Range("L22").Value = 70009
$ws.Range("M22").Value = -15475
$ws.Range("N22").Value = -71067

# GSM!row25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 2698
$ws.Range("I25").Value = 2008
$ws.Range("J25").Value = 2928
$ws.Range("K25").Value = 2008
$ws.Range("L25").Value = 2928
$ws.Range("M25").Value = -1479
$ws.Range("N25").Value = -3986

# GSM!row88
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 29950
$ws.Range("J88").Value = 29950
$ws.Range("L88").Value = 29950
$ws.Range("N88").Value = -30852

# GSM!row91
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H91").Value = 29950
$ws.Range("J91").Value = 29950
$ws.Range("L91").Value = 29950
$ws.Range("N91").Value = -33070

# GSM!row132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3163.2964
$ws.Range("I132").Value = 2713.1875
$ws.Range("K132").Value = 8139.5625
$ws.Range("M132").Value = -5609.5625

# LTW!row40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3359.6
$ws.Range("I40").Value = 1700
$ws.Range("J40").Value = 3774.5
$ws.Range("K40").Value = 1700
$ws.Range("L40").Value = 3774.5
$ws.Range("M40").Value = -1564
$ws.Range("N40").Value = -4046.5

# WVR!row81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 991.8333
$ws.Range("I81").Value = 800.25
$ws.Range("J81").Value = 1087.625
$ws.Range("K81").Value = 1600.5
$ws.Range("L81").Value = 2175.25
$ws.Range("M81").Value = -539.5
$ws.Range("N81").Value = -4297.25

# WVR!row84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 991.8333
$ws.Range("I84").Value = 800.25
$ws.Range("J84").Value = 1087.625
$ws.Range("K84").Value = 8002.5
$ws.Range("L84").Value = 10876.25
$ws.Range("M84").Value = -2698.5
$ws.Range("N84").Value = -21484.25

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 253699.22
$ws.Range("I132").Value = 348826.34
$ws.Range("J132").Value = 41492.54
$ws.Range("K132").Value = 1046479.02
$ws.Range("L132").Value = 124477.62
$ws.Range("M132").Value = -1043949.02
$ws.Range("N132").Value = -129537.62

# WVR!row136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2320.8333
$ws.Range("I136").Value = 1041.6666
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 3124.9998
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -574.9998000000001
$ws.Range("N136").Value = -15900
